$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.040175333333333
$ws.Range("H2").Value = 15.120526
$ws.Range("I2").Value = 0.1638830352839606
$ws.Range("J2").Value = 0.1638830352839606
$ws.Range("Q2").Value = 1.32062674084
$ws.Range("R2").Value = 11.88564066756
$ws.Range("S2").Value = 0.1638830352839606
$ws.Range("T2").Value = 0.1638830352839606

$ws.Range("I3").Value = 0.1966590007046292
$ws.Range("J3").Value = 0.1966590007046292
$ws.Range("S3").Value = 0.1966590007046292
$ws.Range("T3").Value = 0.1966590007046292

$ws.Range("G4").Value = 4.433369666666667
$ws.Range("H4").Value = 13.300109
$ws.Range("I4").Value = 0.1441525402309101
$ws.Range("J4").Value = 0.1441525402309101
$ws.Range("Q4").Value = 1.16163152006
$ws.Range("R4").Value = 10.45468368054
$ws.Range("S4").Value = 0.1441525402309101
$ws.Range("T4").Value = 0.1441525402309101

$ws.Range("G5").Value = 3.607224333333333
$ws.Range("H5").Value = 10.821673
$ws.Range("I5").Value = 0.1172901404415748
$ws.Range("J5").Value = 0.1172901404415748
$ws.Range("Q5").Value = 0.9451649198199998
$ws.Range("R5").Value = 8.506484278379999
$ws.Range("S5").Value = 0.1172901404415748
$ws.Range("T5").Value = 0.1172901404415748

$ws.Range("G6").Value = 6.973136333333334
$ws.Range("H6").Value = 20.919409
$ws.Range("I6").Value = 0.2267339273294199
$ws.Range("J6").Value = 0.2267339273294199
$ws.Range("Q6").Value = 1.82710118206
$ws.Range("R6").Value = 16.44391063854
$ws.Range("S6").Value = 0.2267339273294199
$ws.Range("T6").Value = 0.2267339273294199

$ws.Range("G7").Value = 4.652614333333333
$ws.Range("H7").Value = 13.957843
$ws.Range("I7").Value = 0.1512813560095054
$ws.Range("J7").Value = 0.1512813560095054
$ws.Range("Q7").Value = 1.21907800762
$ws.Range("R7").Value = 10.97170206858
$ws.Range("S7").Value = 0.1512813560095054
$ws.Range("T7").Value = 0.1512813560095054
